# "this added last date 21-11-24"
# Update a handful of quantity/amount figures on the "Product Requisition" sheet
# (Sheet1) and refresh the active view position to where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Data edits -----------------------------------------------------------
# Row 14 (Swap SIM): quantity added
$ws.Range("C14").Value = 75

# Row 31 (lifting qty reduced)
$ws.Range("C31").Value = 10000

# Row 32 (lifting qty added)
$ws.Range("C32").Value = 1000

# Row 43 (DD Payable units adjusted)
$ws.Range("C43").Value = 243909

# --- View state -------------------------------------------------------------
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("K38").Select() | Out-Null
